# Append one new data row (row 66) to each of the four sensor-log sheets,
# mirroring the existing row layout (columns A-I).

$wb = $excel.ActiveWorkbook

$rowsToAdd = @(
    @{
        Sheet = "ROW35-FE-LIFTER"
        A = "2025-03-07 01:42:06"
        B = "0x01,0x90 "
        C = "0x78,0x69,0x90,0x01,0x00,0x00,0x18,0x14,0x41,0x0c,"
        D = "0x01,0x90,"
        E = "0x d"
        F = 400
        G = "'568631262647113770877196"
        H = 400
        I = 13
    },
    @{
        Sheet = "ROW35-MID-LIFTER"
        A = "2025-03-07 01:29:35"
        B = "0x01,0x90 "
        C = "0x78,0x69,0x90,0x01,0x00,0x00,0x18,0x15,0x41,0x0c,"
        D = "0x01,0x90,"
        E = "0x e"
        F = 400
        G = "'568631262647113770942732"
        H = 400
        I = 14
    },
    @{
        Sheet = "ROW02-FE-LIFTER"
        A = "2025-03-07 01:51:45"
        B = "0x01,0x90 "
        C = "0x78,0x69,0x90,0x01,0x00,0x00,0x18,0x06,0x41,0x0c,"
        D = "0x01,0x90,"
        E = "0xff"
        F = 400
        G = "'568631262647113769959692"
        H = 400
        I = 255
    },
    @{
        Sheet = "ROW02-MID-LIFTER"
        A = "2025-03-07 01:41:15"
        B = "0x01,0x90 "
        C = "0x78,0x69,0x90,0x01,0x00,0x00,0x18,0x0b,0x40,0x0c,"
        D = "0x01,0x90,"
        E = "0x 3"
        F = 400
        G = "'568631262647113769959692"
        H = 400
        I = 3
    }
)

foreach ($rowData in $rowsToAdd) {
    $ws = $wb.Worksheets.Item($rowData.Sheet)
    $newRow = 66

    $ws.Range("A$newRow").Value = $rowData.A
    $ws.Range("B$newRow").Value = $rowData.B
    $ws.Range("C$newRow").Value = $rowData.C
    $ws.Range("D$newRow").Value = $rowData.D
    $ws.Range("E$newRow").Value = $rowData.E
    $ws.Range("F$newRow").Value = $rowData.F
    $ws.Range("G$newRow").Value = $rowData.G
    $ws.Range("H$newRow").Value = $rowData.H
    $ws.Range("I$newRow").Value = $rowData.I
}
